# Update the "Förändrad" (Changed) date column (C2:C27) from 2023-12-03
# (serial 45263) to 2023-12-04 (serial 45264) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45263) {
        $cell.Value = 45264
    }
}
